$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.169.82'
$ws.Range('E2').Value = '  +8.45%  '

$ws.Range('D3').Value = '3.316.76'
$ws.Range('E3').Value = '  +4.71%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.17'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '635.55'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.70%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.325'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +18.87%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.611'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.09%  '

$ws.Range('D10').Value = '3.314.25'
$ws.Range('E10').Value = '  +4.91%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.598'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.96%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000272'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +6.03%  '

$ws.Range('E13').Value = '  +2.08%  '

$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.49'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +8.69%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.920.54'
$ws.Range('E15').Value = '  +4.62%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.39'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.94%  '

$ws.Range('D17').Value = '86.888.42'
$ws.Range('E17').Value = '  +8.11%  '

$ws.Range('D18').Value = '3.309.93'
$ws.Range('E18').Value = '  +4.39%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.47'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.46%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.17'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '456.64'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.78%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.03'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.36'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.39'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.40'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +16.06%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.56'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +16.09%  '

$ws.Range('D27').Value = '3.495.72'
$ws.Range('E27').Value = '  +4.81%  '

$ws.Range('B28').Value = 'Litecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '78.24'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.40%  '

$ws.Range('B29').Value = 'Cronos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.217'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +77.76%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000128'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +6.18%  '

$ws.Range('E31').Value = '  -0.23%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.23'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '592.10'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.17%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.11%  '

$ws.Range('E36').Value = '  +3.01%  '

$ws.Range('E37').Value = '  +1.29%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.39'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.61'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +18.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.11%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.417'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.77%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.42'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.16%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.05'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +14.00%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.04'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +13.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '157.82'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.16%  '

$ws.Range('E46').Value = '  +0.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '188.64'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.77'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.51%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.35'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.77%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.784'
$ws.Range('D50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.42'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.34%  '
